$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.440.45"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.916.56"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4713"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").Value = "'0.2855"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "'0.06796"
$ws.Range("E9").Value = "  +4.34%  "
$ws.Range("D10").Value = "'106.66"
$ws.Range("E10").Value = "  +11.71%  "
$ws.Range("D11").Value = "'18.31"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "1.905.96"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "'0.07702"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "'5.205"
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "'0.6579"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "'288.77"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "30.453.93"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "'0.000007637"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'12.94"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "2.151.48"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'5.220"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").Value = "'6.200"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.316"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'168.17"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Value = "'21.53"
$ws.Range("E27").Value = "  +9.63%  "
$ws.Range("D28").Value = "'2.085"
$ws.Range("E28").Value = "  +7.41%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "'1.368"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "'4.163"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "'3.985"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").Value = "'0.05064"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'0.7444"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").Value = "'0.02099"
$ws.Range("E36").Value = "  +9.76%  "
$ws.Range("D37").Value = "'2.742"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").Value = "'2.676"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").Value = "'2.058"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "'109.20"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").Value = "'0.8710"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").Value = "'5.876"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("D43").Value = "'0.4272"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").Value = "'1.0000"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'67.77"
$ws.Range("E45").Value = "  +4.71%  "
$ws.Range("D46").Value = "'50.51"
$ws.Range("E46").Value = "  +18.91%  "
$ws.Range("D47").Value = "'7.186"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").Value = "'9.249"
$ws.Range("E48").Value = "  +3.52%  "
$ws.Range("D49").Value = "'0.1214"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "'34.95"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "'0.3923"
$ws.Range("E51").Value = "  +3.10%  "
